# Regenerate save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# This updates the "K" column (G) values on Sheet1 for rows 2-24.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 2
    3  = 1
    4  = 7
    5  = 6
    6  = 6
    7  = 4
    8  = 7
    9  = 1
    10 = 5
    11 = 1
    12 = 1
    13 = 2
    14 = 3
    15 = 1
    16 = 2
    17 = 1
    18 = 0
    19 = 1
    20 = 1
    21 = 1
    22 = 0
    23 = 2
    24 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
